$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("16").Insert()

$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value = "Los Lagos"
$ws.Range("D16").Value = 44600
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 100112031
$ws.Range("G16").Value = "Poroto verde"
$ws.Range("H16").Value = "Magnum"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 45000
$ws.Range("L16").Value = 45000
$ws.Range("M16").Value = 45000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 1800
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
